$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010991454124451
$ws.Range("B1").Value = 2.126246690750122
$ws.Range("C1").Value = 5.929663181304932
$ws.Range("D1").Value = 1.155157685279846
$ws.Range("E1").Value = 1.185857772827148
